# Scheduled-runner style update of market/profit figures across the
# Sheets workbook. Each worksheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# gets its currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N)
# refreshed for a handful of leve rows with newly-fetched market data.

$wb = $excel.ActiveWorkbook

# NOTE: this interpreter does not bind PowerShell named parameters
# (`-SheetName foo`) reliably, so Set-Row is called positionally below.
function Set-Row {
    param([string]$SheetName, [int]$Row, [hashtable]$Values)
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($col in $Values.Keys) {
        $addr = "$col$Row"
        $ws.Range($addr).Value = $Values[$col]
    }
}

# --- ALC ---
Set-Row "ALC" 55 @{
    H = 873; I = 1194.3334; J = 150; K = 1194.3334; L = 150; M = -980.3334; N = -578
}
Set-Row "ALC" 114 @{
    H = 27562.5; J = 27562.5; L = 27562.5; N = -36240.5
}
Set-Row "ALC" 129 @{
    H = 772650.3; I = 433.1111; J = 950854.25; K = 1299.3333; L = 2852562.75; M = 3700.6667; N = -2862562.75
}
Set-Row "ALC" 138 @{
    H = 3857.2065; I = 2356.3; J = 4274.125; K = 7068.900000000001; L = 12822.375; M = -1928.900000000001; N = -23102.375
}

# --- ARM ---
Set-Row "ARM" 32 @{
    H = 17155.986; I = 13287.111; J = 110009; K = 13287.111; L = 110009; M = -13000.111; N = -110583
}
Set-Row "ARM" 45 @{
    H = 1463.28; I = 1585.6428; J = 1307.5454; K = 1585.6428; L = 1307.5454; M = -1208.6428; N = -2061.5454
}
Set-Row "ARM" 74 @{
    H = 1239.0857; I = 1257.963; J = 1175.375; K = 1257.963; L = 1175.375; M = -383.963; N = -2923.375
}
Set-Row "ARM" 77 @{
    H = 1239.0857; I = 1257.963; J = 1175.375; K = 6289.815; L = 5876.875; M = -1921.815; N = -14612.875
}
Set-Row "ARM" 102 @{
    H = 1703; I = 1781.1111; J = 1000; K = 1781.1111; L = 1000; M = -159.1111000000001; N = -4244
}
Set-Row "ARM" 122 @{
    H = 1650; I = 1680; J = 1500; K = 5040; L = 4500; M = -2590; N = -9400
}

# --- BSM ---
Set-Row "BSM" 99 @{
    H = 33335696; I = 38463744; K = 38463744; M = -38462246
}
Set-Row "BSM" 107 @{
    H = 9951.200000000001; I = 1112.9231; J = 67400; K = 1112.9231; L = 67400; M = 807.0769; N = -71240
}

# --- CRP ---
Set-Row "CRP" 31 @{
    H = 2125.1428; I = 1769.4681; J = 3982.5557; K = 1769.4681; L = 3982.5557; M = -1474.4681; N = -4572.5557
}
Set-Row "CRP" 34 @{
    H = 2125.1428; I = 1769.4681; J = 3982.5557; K = 1769.4681; L = 3982.5557; M = -1567.4681; N = -4386.5557
}
Set-Row "CRP" 134 @{
    H = 55557412; I = 2087.625; J = 500000000; K = 6262.875; L = 1500000000; M = -3727.875; N = -1500005070
}

# --- CUL ---
Set-Row "CUL" 82 @{
    H = 12564.091; I = 2000; J = 13620.5; K = 6000; L = 40861.5; M = -5594; N = -41673.5
}
Set-Row "CUL" 85 @{
    H = 12564.091; I = 2000; J = 13620.5; K = 6000; L = 40861.5; M = -4596; N = -43669.5
}
Set-Row "CUL" 113 @{
    H = 640; I = 955.25; J = 589.5599999999999; K = 2865.75; L = 1768.68; M = -695.75; N = -6108.68
}
Set-Row "CUL" 122 @{
    H = 24012.932; I = 556.8333; J = 27716.525; K = 5011.4997; L = 249448.725; M = -2561.4997; N = -254348.725
}
Set-Row "CUL" 131 @{
    H = 85035.53999999999; J = 73305.92999999999; L = 219917.79; N = -229997.79
}
Set-Row "CUL" 137 @{
    H = 103165; I = 3516.6667; J = 1000000; K = 10550.0001; L = 3000000; M = -5450.000100000001; N = -3010200
}

# --- GSM ---
Set-Row "GSM" 114 @{
    H = 30722; J = 30722; L = 30722; N = -39400
}

# --- LTW ---
Set-Row "LTW" 40 @{
    H = 1124289.1; I = 1124289.1; J = 0; K = 1124289.1; L = 0; M = -1124153.1
}
# N40 no longer exists after this update (folded into M40) - remove it.
$wb.Worksheets.Item("LTW").Range("N40").ClearContents()

Set-Row "LTW" 122 @{
    H = 3012.7827; I = 2832.6365; J = 3177.9167; K = 8497.9095; L = 9533.750100000001; M = -6047.9095; N = -14433.7501
}
Set-Row "LTW" 123 @{
    H = 40000; J = 40000; L = 40000; N = -49800
}

# --- WVR ---
Set-Row "WVR" 122 @{
    H = 1271.6086; I = 1283.9546; J = 1000; K = 3851.8638; L = 3000; M = -1401.8638; N = -7900
}

Write-Output "Applied scheduled Sheets update."
